$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45777
$ws.Range("B2").Value = 2.116666666666667

$ws.Range("A3").Value = 45808
$ws.Range("B3").Value = 9.9
